$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 217-218, pushing the existing rows 217-324
# down to 219-326 (dimension grows from A1:R324 to A1:R326).
$ws.Range("A217:R218").Insert()

# Common (row-invariant) values for this data block.
$marketId = 7
$market   = "Terminal Hortofrutícola Agro Chillán"
$region   = "Ñuble"
$codreg   = 16
$catId    = 100112040
$category = "Cilantro"
$variety  = "Sin especificar"
$unit     = "$/atado 0,5 a 1 kilo"
$kgUnits  = 1
$clasif   = "Hortaliza"

# New row 217: Primera, week of 2023-10-06 (serial 45205)
$r = 217
$ws.Cells.Item($r,1).Value  = $marketId
$ws.Cells.Item($r,2).Value  = $market
$ws.Cells.Item($r,3).Value  = $region
$ws.Cells.Item($r,4).Value  = 45205
$ws.Cells.Item($r,5).Value  = $codreg
$ws.Cells.Item($r,6).Value  = $catId
$ws.Cells.Item($r,7).Value  = $category
$ws.Cells.Item($r,8).Value  = $variety
$ws.Cells.Item($r,9).Value  = "Primera"
$ws.Cells.Item($r,10).Value = 300
$ws.Cells.Item($r,11).Value = 1500
$ws.Cells.Item($r,12).Value = 1500
$ws.Cells.Item($r,13).Value = 1500
$ws.Cells.Item($r,14).Value = $unit
$ws.Cells.Item($r,15).Value = "Provincia de Diguillín"
$ws.Cells.Item($r,16).Value = 1500
$ws.Cells.Item($r,17).Value = $kgUnits
$ws.Cells.Item($r,18).Value = $clasif

# New row 218: Primera, week of 2023-10-06 (serial 45205), Región del Maule
$r = 218
$ws.Cells.Item($r,1).Value  = $marketId
$ws.Cells.Item($r,2).Value  = $market
$ws.Cells.Item($r,3).Value  = $region
$ws.Cells.Item($r,4).Value  = 45205
$ws.Cells.Item($r,5).Value  = $codreg
$ws.Cells.Item($r,6).Value  = $catId
$ws.Cells.Item($r,7).Value  = $category
$ws.Cells.Item($r,8).Value  = $variety
$ws.Cells.Item($r,9).Value  = "Primera"
$ws.Cells.Item($r,10).Value = 300
$ws.Cells.Item($r,11).Value = 2000
$ws.Cells.Item($r,12).Value = 2000
$ws.Cells.Item($r,13).Value = 2000
$ws.Cells.Item($r,14).Value = $unit
$ws.Cells.Item($r,15).Value = "Región del Maule"
$ws.Cells.Item($r,16).Value = 2000
$ws.Cells.Item($r,17).Value = $kgUnits
$ws.Cells.Item($r,18).Value = $clasif
